$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 6495374
$ws.Range("I33").Value = 3707
$ws.Range("J33").Value = 11364124
$ws.Range("K33").Value = 3707
$ws.Range("L33").Value = 11364124
$ws.Range("M33").Value = -3478
$ws.Range("N33").Value = -11364582
# Row 62
$ws.Range("H62").Value = 1944.7407
$ws.Range("I62").Value = 2178.6428
$ws.Range("J62").Value = 1692.8462
$ws.Range("K62").Value = 2178.6428
$ws.Range("L62").Value = 1692.8462
$ws.Range("M62").Value = -1554.6428
$ws.Range("N62").Value = -2940.8462
# Row 65
$ws.Range("H65").Value = 1944.7407
$ws.Range("I65").Value = 2178.6428
$ws.Range("J65").Value = 1692.8462
$ws.Range("K65").Value = 10893.214
$ws.Range("L65").Value = 8464.231
$ws.Range("M65").Value = -7773.214
$ws.Range("N65").Value = -14704.231
# Row 103
$ws.Range("H103").Value = 283532.4
$ws.Range("I103").Value = 916
$ws.Range("J103").Value = 325924.84
$ws.Range("K103").Value = 2748
$ws.Range("L103").Value = 977774.52
$ws.Range("M103").Value = -2162
$ws.Range("N103").Value = -978946.52
# Row 107
$ws.Range("H107").Value = 15625532
$ws.Range("I107").Value = 17857608
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 17857608
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = -17855688
# Row 111
$ws.Range("H111").Value = 125761.5
$ws.Range("I111").Value = 605.4
$ws.Range("J111").Value = 334355
$ws.Range("K111").Value = 1816.2
$ws.Range("L111").Value = 1003065
$ws.Range("M111").Value = 1250.8
$ws.Range("N111").Value = -1009199
# Row 138
$ws.Range("H138").Value = 2384.1516
$ws.Range("I138").Value = 1839.1904
$ws.Range("J138").Value = 3337.8333
$ws.Range("K138").Value = 5517.5712
$ws.Range("L138").Value = 10013.4999
$ws.Range("M138").Value = -377.5712000000003
$ws.Range("N138").Value = -20293.4999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 40
$ws.Range("H40").Value = 1500
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1324
$ws.Range("N40").ClearContents()
# Row 45
$ws.Range("H45").Value = 5449.5386
$ws.Range("I45").Value = 6044.909
$ws.Range("J45").Value = 2175
$ws.Range("K45").Value = 6044.909
$ws.Range("L45").Value = 2175
$ws.Range("M45").Value = -5667.909
$ws.Range("N45").Value = -2929
# Row 122
$ws.Range("H122").Value = 989285.75
$ws.Range("I122").Value = 1224225.2
$ws.Range("J122").Value = 2540
$ws.Range("K122").Value = 3672675.6
$ws.Range("L122").Value = 7620
$ws.Range("M122").Value = -3670225.6
# Row 132
$ws.Range("H132").Value = 3525.4082
$ws.Range("I132").Value = 2658.5
$ws.Range("J132").Value = 4681.2856
$ws.Range("K132").Value = 7975.5
$ws.Range("L132").Value = 14043.8568
$ws.Range("M132").Value = -5445.5
$ws.Range("N132").Value = -19103.8568

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 12772.818
$ws.Range("I20").Value = 1543.1177
$ws.Range("J20").Value = 50953.8
$ws.Range("K20").Value = 1543.1177
$ws.Range("L20").Value = 50953.8
$ws.Range("M20").Value = -1296.1177
# Row 94
$ws.Range("H94").Value = 1793.3529
$ws.Range("I94").Value = 1188.7
$ws.Range("J94").Value = 2657.1428
$ws.Range("K94").Value = 1188.7
$ws.Range("L94").Value = 2657.1428
$ws.Range("M94").Value = -737.7
$ws.Range("N94").Value = -3559.1428
# Row 99
$ws.Range("H99").Value = 58824996
$ws.Range("I99").Value = 90909840
$ws.Range("J99").Value = 2783.3333
$ws.Range("K99").Value = 90909840
$ws.Range("L99").Value = 2783.3333
$ws.Range("M99").Value = -90908342
# Row 107
$ws.Range("H107").Value = 637.0526
$ws.Range("I107").Value = 891.2
$ws.Range("J107").Value = 354.66666
$ws.Range("K107").Value = 891.2
$ws.Range("L107").Value = 354.66666
$ws.Range("M107").Value = 1028.8
$ws.Range("N107").Value = -4194.66666
# Row 134
$ws.Range("H134").Value = 38796.434
$ws.Range("I134").Value = 6848.95
$ws.Range("J134").Value = 102691.4
$ws.Range("K134").Value = 20546.85
$ws.Range("L134").Value = 308074.2
$ws.Range("M134").Value = -18011.85
$ws.Range("N134").Value = -313144.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1742.3158
$ws.Range("I16").Value = 1263.875
$ws.Range("J16").Value = 2090.2727
$ws.Range("K16").Value = 1263.875
$ws.Range("L16").Value = 2090.2727
$ws.Range("M16").Value = -976.875
# Row 105
$ws.Range("H105").Value = 2573.36
$ws.Range("I105").Value = 2520.1667
$ws.Range("J105").Value = 3850
$ws.Range("K105").Value = 2520.1667
$ws.Range("L105").Value = 3850
$ws.Range("M105").Value = -773.1667000000002
# Row 113
$ws.Range("H113").Value = 1742.3158
$ws.Range("I113").Value = 1263.875
$ws.Range("J113").Value = 2090.2727
$ws.Range("K113").Value = 1263.875
$ws.Range("L113").Value = 2090.2727
$ws.Range("M113").Value = 906.125

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 61
$ws.Range("H61").Value = 525
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 525
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 1575
$ws.Range("N61").Value = -2005
$ws.Range("M61").ClearContents()
# Row 113
$ws.Range("H113").Value = 1000597.25
$ws.Range("I113").Value = 1000605.8
$ws.Range("J113").Value = 1000554.5
$ws.Range("K113").Value = 3001817.4
$ws.Range("L113").Value = 3001663.5
$ws.Range("M113").Value = -2999647.4
$ws.Range("N113").Value = -3006003.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 865.5
$ws.Range("I31").Value = 865.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 865.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -573.5
# Row 37
$ws.Range("H37").Value = 865.5
$ws.Range("I37").Value = 865.5
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 865.5
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -588.5
# Row 97
$ws.Range("H97").Value = 1221.8846
$ws.Range("I97").Value = 1221.8846
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1221.8846
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -725.8846000000001
# Row 113
$ws.Range("H113").Value = 40001456
$ws.Range("I113").Value = 90910090
$ws.Range("J113").Value = 1815.1428
$ws.Range("K113").Value = 90910090
$ws.Range("L113").Value = 1815.1428
$ws.Range("M113").Value = -90907920
$ws.Range("N113").Value = -6155.1428
# Row 122
$ws.Range("H122").Value = 38672564
$ws.Range("I122").Value = 42594004
$ws.Range("J122").Value = 27779678
$ws.Range("K122").Value = 127782012
$ws.Range("L122").Value = 83339034
$ws.Range("M122").Value = -127779562
# Row 123
$ws.Range("H123").Value = 9480.296
$ws.Range("I123").Value = 9000
$ws.Range("J123").Value = 9498.77
$ws.Range("K123").Value = 9000
$ws.Range("L123").Value = 9498.77
$ws.Range("M123").Value = -6550
$ws.Range("N123").Value = -14398.77
# Row 126
$ws.Range("H126").Value = 6140.92
$ws.Range("I126").Value = 8145.25
$ws.Range("J126").Value = 2577.6667
$ws.Range("K126").Value = 24435.75
$ws.Range("L126").Value = 7733.000100000001
$ws.Range("M126").Value = -21965.75
$ws.Range("N126").Value = -12673.0001
# Row 132
$ws.Range("H132").Value = 3345.359
$ws.Range("I132").Value = 4058.6843
$ws.Range("J132").Value = 2667.7
$ws.Range("K132").Value = 12176.0529
$ws.Range("L132").Value = 8003.099999999999
$ws.Range("M132").Value = -9646.052899999999
$ws.Range("N132").Value = -13063.1

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1292.8125
$ws.Range("I16").Value = 1292.8125
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1292.8125
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1122.8125
$ws.Range("N16").ClearContents()
# Row 40
$ws.Range("H40").Value = 25001834
$ws.Range("I40").Value = 32259924
$ws.Range("J40").Value = 1745
$ws.Range("K40").Value = 32259924
$ws.Range("L40").Value = 1745
$ws.Range("M40").Value = -32259788
$ws.Range("N40").Value = -2017
# Row 61
$ws.Range("H61").Value = 1647.1364
$ws.Range("I61").Value = 1338.3572
$ws.Range("J61").Value = 2187.5
$ws.Range("K61").Value = 1338.3572
$ws.Range("L61").Value = 2187.5
$ws.Range("M61").Value = -1136.3572
$ws.Range("N61").Value = -2591.5
# Row 113
$ws.Range("H113").Value = 1647.1364
$ws.Range("I113").Value = 1338.3572
$ws.Range("J113").Value = 2187.5
$ws.Range("K113").Value = 1338.3572
$ws.Range("L113").Value = 2187.5
$ws.Range("M113").Value = 831.6428000000001
$ws.Range("N113").Value = -6527.5
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
# Row 122
$ws.Range("H122").Value = 3259858.2
$ws.Range("I122").Value = 3762046
$ws.Range("J122").Value = 1669596.5
$ws.Range("K122").Value = 11286138
$ws.Range("L122").Value = 5008789.5
$ws.Range("M122").Value = -11283688
$ws.Range("N122").Value = -5013689.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 166667360
$ws.Range("I107").Value = 200000580
$ws.Range("J107").Value = 1280
$ws.Range("K107").Value = 600001740
$ws.Range("L107").Value = 3840
$ws.Range("M107").Value = -599999820
# Row 122
$ws.Range("H122").Value = 2182.7646
$ws.Range("I122").Value = 1629.7
$ws.Range("J122").Value = 2972.8572
$ws.Range("K122").Value = 4889.1
$ws.Range("L122").Value = 8918.571599999999
$ws.Range("M122").Value = -2439.1
$ws.Range("N122").Value = -13818.5716
